$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains its text formatting so values such as
# "1.00" or "0.660" do not get auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.163.21"
$ws.Range("E2").Value = "  -3.04%  "
$ws.Range("D3").Value = "2.953.74"
$ws.Range("E3").Value = "  -2.44%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "519.18"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "128.22"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "2.946.48"
$ws.Range("E8").Value = "  -2.54%  "
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("D10").Value = "6.06"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").Value = "0.145"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("E12").Value = "  -2.67%  "
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "32.67"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").Value = "3.439.99"
$ws.Range("E15").Value = "  -2.13%  "
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "60.289.58"
$ws.Range("E17").Value = "  -2.72%  "
$ws.Range("D18").Value = "2.959.52"
$ws.Range("E18").Value = "  -2.09%  "
$ws.Range("D19").Value = "6.39"
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("D20").Value = "450.71"
$ws.Range("E20").Value = "  -4.58%  "
$ws.Range("D21").Value = "12.83"
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D22").Value = "0.660"
$ws.Range("E22").Value = "  -3.29%  "
$ws.Range("D23").Value = "6.71"
$ws.Range("E23").Value = "  -2.47%  "
$ws.Range("D24").Value = "77.22"
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("D25").Value = "11.60"
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "2.60"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").Value = "7.59"
$ws.Range("E28").Value = "  -4.73%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("E30").Value = "  +3.46%  "
$ws.Range("D31").Value = "24.79"
$ws.Range("E31").Value = "  -2.21%  "
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("D33").Value = "54.41"
$ws.Range("E33").Value = "  -3.77%  "
$ws.Range("D34").Value = "2.21"
$ws.Range("E34").Value = "  -4.49%  "
$ws.Range("D35").Value = "5.28"
$ws.Range("E35").Value = "  +1.98%  "
$ws.Range("D36").Value = "5.66"
$ws.Range("E36").Value = "  -2.40%  "
$ws.Range("D37").Value = "445.36"
$ws.Range("E37").Value = "  -4.28%  "
$ws.Range("D38").Value = "3.142.34"
$ws.Range("E38").Value = "  +2.75%  "
$ws.Range("D39").Value = "0.0766"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("D40").Value = "0.0372"
$ws.Range("E40").Value = "  -2.82%  "
$ws.Range("D41").Value = "0.114"
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("D42").Value = "7.88"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("E43").Value = "  -4.74%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").Value = "0.239"
$ws.Range("E45").Value = "  -1.76%  "
$ws.Range("D46").Value = "24.92"
$ws.Range("E46").Value = "  +4.73%  "
$ws.Range("D47").Value = "116.98"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "1.90"
$ws.Range("E49").Value = "  -3.73%  "
$ws.Range("D50").Value = "0.0₃0498"
$ws.Range("E50").Value = "  -5.99%  "
$ws.Range("E51").Value = "  +8.28%  "
